# Trade #61 closed at 2026-02-17 08:48:23 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Summary" sheet - update the rollup metrics now that trade #61 closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.69   # Current Capital
$summary.Range("B4").Value = -0.31    # Total P&L $
$summary.Range("B5").Value = -0.1     # Total P&L %
$summary.Range("B6").Value = 61        # Total Trades
$summary.Range("B7").Value = 25        # Winning Trades
$summary.Range("B9").Value = 40.98     # Win Rate %

# ---------------------------------------------------------------------------
# 2. "Strategy Status" sheet - update the MarketMaking strategy row (row 4).
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.69      # Capital
$status.Range("D4").Value = 61         # Trades
$status.Range("E4").Value = -0.31     # P&L $
$status.Range("F4").Value = -0.31     # P&L %
$status.Range("G4").Value = 40.98      # Win Rate %

# ---------------------------------------------------------------------------
# 3. Append the newly-closed trade (#61) to both the "All Trades" and the
#    "MarketMaking" trade logs - they mirror each other row for row.
# ---------------------------------------------------------------------------
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 62

    $ws.Cells.Item($row, 1).Value = 61

    # Date column would otherwise be auto-parsed into a date serial by
    # Excel's smart-entry; force it to stay plain text like the rest of
    # the column, then drop the leftover text-format style so the cell
    # is indistinguishable from a normal, unstyled cell.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "08:48:16"
    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.78
    $ws.Cells.Item($row, 7).Value = 0.79
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 1.2821
    $ws.Cells.Item($row, 10).Value = 0.01
    $ws.Cells.Item($row, 11).Value = 99.69
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.14
}
